$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 7 & 8 (selfemployed / unemployed volumes): caveat text now refers to
# "all ages (16+)" instead of "16-64 year olds".
$caveatGeneric = @"
<ol>
  <li>Figures are for all ages (16+).</li>
<li>Each estimate from the Annual Population Survey carries a margin of error. These are available in the original data via NOMIS. Large margins of error are usually associated with groups with only a small number of respondents. Therefore, please take caution when interpreting data from small subgroups.</li>
<li>Use caution when interpreting this data. A difference between subgroups does not necessarily imply any causality. There could be other contributing factors at work.</li>
</ol>
"@
$ws.Range("F7").Value = $caveatGeneric
$ws.Range("F8").Value = $caveatGeneric

# Row 4 (unemployedRate): caveat text now refers to "all ages (16+)" and
# includes an extra clarifying bullet about the unemployment rate definition.
$caveatRow4 = @"
<ol>
  <li>Figures are for all ages (16+).</li>
<li>The unemployment rate is not the proportion of the total population that is unemployed. It is the proportion of the economically active population (that is, those in work plus those seeking and available to work) that is unemployed.</li>
<li>Each estimate from the Annual Population Survey carries a margin of error. These are available in the original data via NOMIS. Large margins of error are usually associated with groups with only a small number of respondents. Therefore, please take caution when interpreting data from small subgroups.</li>
<li>Use caution when interpreting this data. A difference between subgroups does not necessarily imply any causality. There could be other contributing factors at work.</li>
</ol>
"@
$ws.Range("F4").Value = $caveatRow4

# Row 6 (inemployment / employment volumes): caveat text now refers to
# "all ages (16+)" instead of "16-64 year olds", keeping the
# SOC2020/SIC2007 classification bullets.
$caveatRow6 = @"
<ol>
  <li>Figures are for all ages (16+).</li>
<li>Each estimate from the Annual Population Survey carries a margin of error. These are available in the original data via NOMIS. Large margins of error are usually associated with groups with only a small number of respondents. Therefore, please take caution when interpreting data from small subgroups.</li>
<li>Use caution when interpreting this data. A difference between subgroups does not necessarily imply any causality. There could be other contributing factors at work.</li>
<li>Standard Occupational Classification 2020 (SOC2020).</li>
<li>Industry groups are Standard Industrial Classification: SIC 2007.</li>
</ol>
"@
$ws.Range("F6").Value = $caveatRow6

# The longer caveat text in rows 4 and 6 reflows to a taller row.
$ws.Rows.Item(4).RowHeight = 261.5
$ws.Rows.Item(6).RowHeight = 247

# Restore the view to the top of the sheet with A2 selected (matches a
# fresh open/save rather than the mid-scroll state the file was left in).
[void]$ws.Range("A2").Select()
